$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28. This shifts the existing rows 28..40
# down to 29..41 and copies formatting (incl. the date number format on
# column D) from the row above, matching the target workbook.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new record.
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = 'Vega Monumental Concepción'
$ws.Range("C28").Value = 'Bíobío'
$ws.Range("D28").Value = 45007
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = 'Arveja Verde'
$ws.Range("H28").Value = 'Perfection'
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 27000
$ws.Range("L28").Value = 28000
$ws.Range("M28").Value = 27500
$ws.Range("N28").Value = '$/saco 25 kilos'
$ws.Range("O28").Value = 'Región Metropolitana'
$ws.Range("P28").Value = 1100
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = 'Hortaliza'
